$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1777.6666
$ws.Range("J17").Value = 1777.6666
$ws.Range("L17").Value = 5332.9998
$ws.Range("N17").Value = -5668.9998
$ws.Range("H140").Value = 83956
$ws.Range("J140").Value = 83956
$ws.Range("L140").Value = 83956
$ws.Range("N140").Value = -94316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1390.5
$ws.Range("I2").Value = 1217.3
$ws.Range("K2").Value = 1217.3
$ws.Range("M2").Value = -1104.3
$ws.Range("H32").Value = 4821.5
$ws.Range("I32").Value = 2335.1428
$ws.Range("J32").Value = 17874.875
$ws.Range("K32").Value = 2335.1428
$ws.Range("L32").Value = 17874.875
$ws.Range("M32").Value = -2048.1428
$ws.Range("N32").Value = -18448.875
$ws.Range("H74").Value = 8073059.5
$ws.Range("I74").Value = 12502922
$ws.Range("J74").Value = 18764.363
$ws.Range("K74").Value = 12502922
$ws.Range("L74").Value = 18764.363
$ws.Range("M74").Value = -12502048
$ws.Range("N74").Value = -20512.363
$ws.Range("H77").Value = 8073059.5
$ws.Range("I77").Value = 12502922
$ws.Range("J77").Value = 18764.363
$ws.Range("K77").Value = 62514610
$ws.Range("L77").Value = 93821.815
$ws.Range("M77").Value = -62510242
$ws.Range("N77").Value = -102557.815
$ws.Range("H102").Value = 3403.625
$ws.Range("I102").Value = 3982.5
$ws.Range("J102").Value = 509.25
$ws.Range("K102").Value = 3982.5
$ws.Range("L102").Value = 509.25
$ws.Range("M102").Value = -2360.5
$ws.Range("N102").Value = -3753.25
$ws.Range("H116").Value = 1390.5
$ws.Range("I116").Value = 1217.3
$ws.Range("K116").Value = 1217.3
$ws.Range("M116").Value = 1076.7
$ws.Range("H122").Value = 2233.608
$ws.Range("I122").Value = 1700.4
$ws.Range("K122").Value = 5101.200000000001
$ws.Range("M122").Value = -2651.200000000001
$ws.Range("H132").Value = 7756270
$ws.Range("I132").Value = 10754927
$ws.Range("K132").Value = 32264781
$ws.Range("M132").Value = -32262251

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1390.5
$ws.Range("I3").Value = 1217.3
$ws.Range("K3").Value = 1217.3
$ws.Range("M3").Value = -1103.3
$ws.Range("H20").Value = 5928.778
$ws.Range("I20").Value = 5032.9443
$ws.Range("J20").Value = 7720.4443
$ws.Range("K20").Value = 5032.9443
$ws.Range("L20").Value = 7720.4443
$ws.Range("M20").Value = -4785.9443
$ws.Range("N20").Value = -8214.444299999999
$ws.Range("H94").Value = 1615.5
$ws.Range("I94").Value = 656.3570999999999
$ws.Range("J94").Value = 2958.3
$ws.Range("K94").Value = 656.3570999999999
$ws.Range("L94").Value = 2958.3
$ws.Range("M94").Value = -205.3570999999999
$ws.Range("N94").Value = -3860.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 595273.25
$ws.Range("I31").Value = 10777.375
$ws.Range("J31").Value = 984937.2
$ws.Range("K31").Value = 10777.375
$ws.Range("L31").Value = 984937.2
$ws.Range("M31").Value = -10482.375
$ws.Range("N31").Value = -985527.2
$ws.Range("H34").Value = 595273.25
$ws.Range("I34").Value = 10777.375
$ws.Range("J34").Value = 984937.2
$ws.Range("K34").Value = 10777.375
$ws.Range("L34").Value = 984937.2
$ws.Range("M34").Value = -10575.375
$ws.Range("N34").Value = -985341.2
$ws.Range("H99").Value = 2496.4285
$ws.Range("I99").Value = 2200.1
$ws.Range("K99").Value = 2200.1
$ws.Range("M99").Value = -702.0999999999999
$ws.Range("H105").Value = 1952.5454
$ws.Range("I105").Value = 2158.5
$ws.Range("J105").Value = 1403.3334
$ws.Range("K105").Value = 2158.5
$ws.Range("L105").Value = 1403.3334
$ws.Range("M105").Value = -411.5
$ws.Range("N105").Value = -4897.3334
$ws.Range("H126").Value = 2496.4285
$ws.Range("I126").Value = 2200.1
$ws.Range("K126").Value = 6600.299999999999
$ws.Range("M126").Value = -4130.299999999999
$ws.Range("H132").Value = 9173.625
$ws.Range("I132").Value = 4727.4375
$ws.Range("K132").Value = 14182.3125
$ws.Range("M132").Value = -11652.3125
$ws.Range("H134").Value = 2561.9119
$ws.Range("I134").Value = 1605.5
$ws.Range("J134").Value = 3928.2144
$ws.Range("K134").Value = 4816.5
$ws.Range("L134").Value = 11784.6432
$ws.Range("M134").Value = -2281.5
$ws.Range("N134").Value = -16854.6432
$ws.Range("H141").Value = 220852.33
$ws.Range("I141").Value = 39607.668
$ws.Range("K141").Value = 39607.668
$ws.Range("M141").Value = -34427.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1307.5238
$ws.Range("J122").Value = 1436.5555
$ws.Range("L122").Value = 12928.9995
$ws.Range("N122").Value = -17828.9995
$ws.Range("H132").Value = 3179.8
$ws.Range("J132").Value = 3916.5
$ws.Range("L132").Value = 35248.5
$ws.Range("N132").Value = -40308.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H43").Value = 28333
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 32999.6
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 32999.6
$ws.Range("M43").Value = -4849
$ws.Range("N43").Value = -33301.6
$ws.Range("H80").Value = 5621.25
$ws.Range("I80").Value = 4996.8335
$ws.Range("J80").Value = 7494.5
$ws.Range("K80").Value = 4996.8335
$ws.Range("L80").Value = 7494.5
$ws.Range("M80").Value = -3998.8335
$ws.Range("N80").Value = -9490.5
$ws.Range("H83").Value = 5621.25
$ws.Range("I83").Value = 4996.8335
$ws.Range("J83").Value = 7494.5
$ws.Range("K83").Value = 24984.1675
$ws.Range("L83").Value = 37472.5
$ws.Range("M83").Value = -19992.1675
$ws.Range("N83").Value = -47456.5
$ws.Range("H97").Value = 1349.174
$ws.Range("I97").Value = 776.3125
$ws.Range("K97").Value = 776.3125
$ws.Range("M97").Value = -280.3125
$ws.Range("H113").Value = 3729.0344
$ws.Range("J113").Value = 4101.6113
$ws.Range("L113").Value = 4101.6113
$ws.Range("N113").Value = -8441.6113
$ws.Range("H132").Value = 29415262
$ws.Range("I132").Value = 32261482
$ws.Range("K132").Value = 96784446
$ws.Range("M132").Value = -96781916
$ws.Range("H136").Value = 9188.842000000001
$ws.Range("J136").Value = 9188.842000000001
$ws.Range("L136").Value = 27566.526
$ws.Range("N136").Value = -32666.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1395.375
$ws.Range("I22").Value = 1277.3334
$ws.Range("K22").Value = 1277.3334
$ws.Range("M22").Value = -982.3334
$ws.Range("H27").Value = 1395.375
$ws.Range("I27").Value = 1277.3334
$ws.Range("K27").Value = 1277.3334
$ws.Range("M27").Value = -1170.3334
$ws.Range("H46").Value = 3689.2222
$ws.Range("I46").Value = 3160.182
$ws.Range("J46").Value = 4520.5713
$ws.Range("K46").Value = 3160.182
$ws.Range("L46").Value = 4520.5713
$ws.Range("M46").Value = -2972.182
$ws.Range("N46").Value = -4896.5713
$ws.Range("H68").Value = 2552.7778
$ws.Range("I68").Value = 2425
$ws.Range("K68").Value = 2425
$ws.Range("M68").Value = -1676
$ws.Range("H71").Value = 2552.7778
$ws.Range("I71").Value = 2425
$ws.Range("K71").Value = 12125
$ws.Range("M71").Value = -8381
$ws.Range("H136").Value = 52210.035
$ws.Range("I136").Value = 11395.883
$ws.Range("K136").Value = 34187.649
$ws.Range("M136").Value = -31637.649

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15396746
$ws.Range("I62").Value = 2949.5
$ws.Range("K62").Value = 2949.5
$ws.Range("M62").Value = -2325.5
$ws.Range("H65").Value = 15396746
$ws.Range("I65").Value = 2949.5
$ws.Range("K65").Value = 14747.5
$ws.Range("M65").Value = -11627.5
$ws.Range("H75").Value = 75000000
$ws.Range("J75").Value = 75000000
$ws.Range("L75").Value = 75000000
$ws.Range("N75").Value = -75001872
$ws.Range("H78").Value = 75000000
$ws.Range("J78").Value = 75000000
$ws.Range("L78").Value = 225000000
$ws.Range("N78").Value = -225009360
$ws.Range("H100").Value = 1075.1904
$ws.Range("I100").Value = 1153.2
$ws.Range("K100").Value = 2306.4
$ws.Range("M100").Value = -1765.4
$ws.Range("H132").Value = 318088.38
$ws.Range("I132").Value = 4820.269
$ws.Range("K132").Value = 14460.807
$ws.Range("M132").Value = -11930.807
